$d = $word.ActiveDocument

function New-OpenXmlPackage($bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Insert a new empty (bold / orange / 32pt) paragraph right before the
#    last paragraph of the document (the one that currently only contains "-").
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$dashPara = $d.Paragraphs.Item($lastIndex)
$beforeDash = $d.Range($dashPara.Range.Start, $dashPara.Range.Start)

$boldParaBody = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="ED7D31" w:themeColor="accent2"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-MX"/></w:rPr></w:pPr></w:p>'
$beforeDash.InsertXML((New-OpenXmlPackage $boldParaBody))

# ---------------------------------------------------------------------------
# 2) Append the new sentence onto the existing "-" paragraph (now the last
#    paragraph again, since the insert above only added a paragraph before it).
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$dashPara = $d.Paragraphs.Item($lastIndex)
$dashPara.Range.InsertAfter(" Al final del proyecto no implementamos algunas inferencias que propusimos al principio como diseño de la solución. Como la de existeGanador() o las encargadas de llevar un conteo de puntos. Esto se debe a que realmente no lo vimos necesario por como se desarrolló la interfaz gráfica. Además, un sistema de puntos estaba de más y no aportaba mucho sentido a nuestra idea de juego como tal. También se corrigieron algunas o se les cambió el nombre. ")

# ---------------------------------------------------------------------------
# 3) Append the remaining new paragraphs (with blank separators) after the
#    "-" paragraph, at the very end of the document body.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$dashPara = $d.Paragraphs.Item($lastIndex)
$endOfDoc = $d.Range($dashPara.Range.End, $dashPara.Range.End)

$rPr28 = '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-MX"/></w:rPr>'
$emptyPara = '<w:p><w:pPr>' + $rPr28 + '</w:pPr></w:p>'

$p1 = '<w:p><w:pPr>' + $rPr28 + '</w:pPr>' `
    + '<w:r>' + $rPr28 + '<w:t>-Para cumplir con los procedimientos dinámicos se aplicó un poco más el uso de retract y assertz a algunas inferencias como a la de tirarDado(), generarPistaRandom()</w:t></w:r>' `
    + '<w:r>' + $rPr28 + '<w:t>,quitarPistaRandom(), etc.</w:t></w:r>' `
    + '</w:p>'

$p2 = '<w:p><w:pPr>' + $rPr28 + '</w:pPr>' `
    + '<w:r>' + $rPr28 + '<w:t xml:space="preserve">-Hubo complicaciones con el uso de JavaFX. Ya que uno de los integrantes del grupo tenía conocimientos </w:t></w:r>' `
    + '<w:r>' + $rPr28 + '<w:t>con Java normal manipulando JFrames</w:t></w:r>' `
    + '<w:r>' + $rPr28 + '<w:t>, nos inclinamos por crear la interfaz gráfica de esa manera.</w:t></w:r>' `
    + '</w:p>'

$p3 = '<w:p><w:pPr>' + $rPr28 + '</w:pPr>' `
    + '<w:r>' + $rPr28 + '<w:t>-A pesar de que en un principio nos parecía extraño juntar prolog con java, al final nos pareció una experiencia muy chiva. Crear un juego en base a resultados ‘true’ o ‘false’ es algo que no habíamos enfrentado antes.</w:t></w:r>' `
    + '</w:p>'

$tailBody = $emptyPara + $p1 + $emptyPara + $p2 + $emptyPara + $p3
$endOfDoc.InsertXML((New-OpenXmlPackage $tailBody))

Write-Host "Final paragraph count:" $d.Paragraphs.Count
